$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DQ_Metrics")

# Insert a new column at J, shifting old J:T (tracerCase_rel_py_ipat .. tracerCase_no_py)
# one column to the right, to K:U.
$ws.Columns("J:J").Insert()

# Rename the metric headers whose meaning changed (same column position).
$ws.Range("C1").Value = "item_completeness_rate"
$ws.Range("D1").Value = "value_completeness_rate"
$ws.Range("F1").Value = "range_plausibility_rate"
$ws.Range("I1").Value = "rdCase_dissimilarity_rate"

# Populate the newly inserted column header.
$ws.Range("J1").Value = "rdCase_rel_py_ipat"

# Append the new trailing headers.
$ws.Range("V1").Value = "missing_item_no_py"
$ws.Range("W1").Value = "missing_value_no_py"
$ws.Range("X1").Value = "orphaMissing_no_py"
$ws.Range("Y1").Value = "implausible_codeLink_no_py"
$ws.Range("Z1").Value = "outlier_no_py"
$ws.Range("AA1").Value = "ambigous_rdCase_no_py"
$ws.Range("AB1").Value = "duplicateRdCase_no_py"

# Update the data row values that changed.
$ws.Range("C2").Value = 100
$ws.Range("D2").Value = 59.51
$ws.Range("F2").Value = 93.75
$ws.Range("I2").Value = 94.12
$ws.Range("J2").Value = 0.17

# Populate the new trailing data values.
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = 149
$ws.Range("X2").Value = 2
$ws.Range("Y2").Value = 10
$ws.Range("Z2").Value = 5
$ws.Range("AA2").Value = 10
$ws.Range("AB2").Value = 1
